# Legislator property workbook edit
# Sheet 6 = "基金受益憑證" (fund / beneficiary certificates).
#
# Today the sheet has no real header row - row 1 is just a duplicate of
# row 2's data - and rows 2-3 are missing the property_category /
# category / date / legislator_name / legislator_id / source_file /
# index metadata columns that every other sheet in this workbook already
# carries (see sheet1 "土地" / sheet5 "股票" for the pattern). This
# brings sheet6 in line with the rest of the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(6)

# Helper: write $text into $targetCell as a literal text value, even when
# it looks like a date (e.g. "2012-03-26") which Excel would otherwise
# silently coerce into a date serial number. Goes through a scratch cell
# as a formula so the result lands as a plain string value with no
# number-format/style side effects.
function Set-LiteralText {
    param($targetCell, [string]$text)
    $scratch = $ws.Range("Z100")
    $scratch.Formula = '="' + $text + '"'
    $scratch.Copy()
    $targetCell.PasteSpecial(-4163) # xlPasteValues
    $scratch.Clear()
    $excel.CutCopyMode = $false
}

# --- Row 1: turn the duplicated data row into a proper header row ------
$ws.Cells.Item(1, 2).Value = "name"
$ws.Cells.Item(1, 3).Value = "owner"
$ws.Cells.Item(1, 4).Value = "dealer"
$ws.Cells.Item(1, 5).Value = "quantity"
$ws.Cells.Item(1, 6).Value = "face_value"
$ws.Cells.Item(1, 7).Value = "currency"
$ws.Cells.Item(1, 8).Value = "total"

# New header cells I1:O1 - copy the existing header style (bold + border,
# same as B1:H1) onto them before filling in the labels.
$ws.Range("B1").Copy()
$ws.Range("I1:O1").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Cells.Item(1, 9).Value = "property_category"
$ws.Cells.Item(1, 10).Value = "category"
$ws.Cells.Item(1, 11).Value = "date"
$ws.Cells.Item(1, 12).Value = "legislator_name"
$ws.Cells.Item(1, 13).Value = "legislator_id"
$ws.Cells.Item(1, 14).Value = "source_file"
$ws.Cells.Item(1, 15).Value = "index"

# --- Row 2 data (existing fund holding, unchanged) ----------------------
$ws.Cells.Item(2, 2).Value = "貝萊德世界礦業"
$ws.Cells.Item(2, 3).Value = "潘維剛"
$ws.Cells.Item(2, 4).Value = "永豐商業銀行"
$ws.Cells.Item(2, 5).Value = 5669.54
$ws.Cells.Item(2, 6).Value = 65.47
$ws.Cells.Item(2, 7).Value = "美金"
$ws.Cells.Item(2, 8).Value = 10975934.06

$ws.Cells.Item(2, 9).Value = "fund"
$ws.Cells.Item(2, 10).Value = "normal"
Set-LiteralText $ws.Cells.Item(2, 11) "2012-03-26"
$ws.Cells.Item(2, 12).Value = "潘維剛"
$ws.Cells.Item(2, 13).Value = 678
$ws.Cells.Item(2, 14).Value = "tmp71a01"
$ws.Cells.Item(2, 15).Value = 81

# --- Row 3 data (existing fund holding, unchanged) ----------------------
$ws.Cells.Item(3, 2).Value = "聯博全球高收益債券AT股"
$ws.Cells.Item(3, 3).Value = "潘維剛"
$ws.Cells.Item(3, 4).Value = "永豐商業銀行"
$ws.Cells.Item(3, 5).Value = 44345.898
$ws.Cells.Item(3, 6).Value = 4.51
$ws.Cells.Item(3, 7).Value = "美金"
$ws.Cells.Item(3, 8).Value = 5914000

$ws.Cells.Item(3, 9).Value = "fund"
$ws.Cells.Item(3, 10).Value = "normal"
Set-LiteralText $ws.Cells.Item(3, 11) "2012-03-26"
$ws.Cells.Item(3, 12).Value = "潘維剛"
$ws.Cells.Item(3, 13).Value = 678
$ws.Cells.Item(3, 14).Value = "tmp71a01"
$ws.Cells.Item(3, 15).Value = 82
